$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as sourced from the commit diff.
# Values are written as text strings so the original inlineStr cell formatting/content is preserved exactly.

$ws.Range('D2').Value = '35.598.98'
$ws.Range('E2').Value = '  +1.89%  '
$ws.Range('D3').Value = '1.904.03'
$ws.Range('E3').Value = '  +3.33%  '
$ws.Range('E4').Value = '  +0.48%  '
$ws.Range('D5').Value = '244.95'
$ws.Range('E5').Value = '  +5.31%  '
$ws.Range('D6').Value = '0.635'
$ws.Range('E6').Value = '  +2.32%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').Value = '42.62'
$ws.Range('E8').Value = '  +3.36%  '
$ws.Range('D9').Value = '0.338'
$ws.Range('E9').Value = '  +3.03%  '
$ws.Range('D10').Value = '0.0707'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').Value = '2.176.34'
$ws.Range('D13').Value = '12.54'
$ws.Range('E13').Value = '  +10.29%  '
$ws.Range('D14').Value = '1.903.53'
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('D15').Value = '0.694'
$ws.Range('E15').Value = '  +3.47%  '
$ws.Range('D16').Value = '4.84'
$ws.Range('E16').Value = '  +3.78%  '
$ws.Range('D17').Value = '35.573.55'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').Value = '72.16'
$ws.Range('E18').Value = '  +3.23%  '
$ws.Range('D19').Value = '0.0₃0810'
$ws.Range('E19').Value = '  +2.61%  '
$ws.Range('D20').Value = '244.71'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').Value = '12.47'
$ws.Range('E21').Value = '  +2.65%  '
$ws.Range('D22').Value = '4.93'
$ws.Range('E22').Value = '  +3.41%  '
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').Value = '2.28'
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('D25').Value = '171.09'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('D26').Value = '2.11'
$ws.Range('E26').Value = '  +28.76%  '
$ws.Range('D27').Value = '8.50'
$ws.Range('E27').Value = '  +8.81%  '
$ws.Range('D28').Value = '17.98'
$ws.Range('E28').Value = '  +3.07%  '
$ws.Range('E29').Value = '  +1.73%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '0.955'
$ws.Range('E30').Value = '  +28.70%  '
$ws.Range('D31').Value = '4.10'
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0567'
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('E34').Value = '  +5.65%  '
$ws.Range('E35').Value = '  +8.45%  '
$ws.Range('D36').Value = '2.05'
$ws.Range('E36').Value = '  +5.13%  '
$ws.Range('E37').Value = '  +8.40%  '
$ws.Range('E38').Value = '  +3.88%  '
$ws.Range('D39').Value = '0.0206'
$ws.Range('E39').Value = '  +5.02%  '
$ws.Range('D40').Value = '91.16'
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D41').Value = '1.359.08'
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('D42').Value = '15.25'
$ws.Range('E42').Value = '  +5.01%  '
$ws.Range('D43').Value = '0.0601'
$ws.Range('E43').Value = '  +13.40%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '2.36'
$ws.Range('E44').Value = '  +4.15%  '
$ws.Range('B45').Value = 'Gas'
$ws.Range('C45').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D45').Value = '12.98'
$ws.Range('E45').Value = '  +31.47%  '
$ws.Range('D46').Value = '46.26'
$ws.Range('E46').Value = '  +36.37%  '
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').Value = '6.68'
$ws.Range('E49').Value = '  +5.73%  '
$ws.Range('D50').Value = '2.084.61'
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('E51').Value = '  +3.41%  '
